$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1): F2 357 -> 358, F4 76 -> 78
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 358
$ws1.Range("F4").Value = 78

# Update "全部类型" sheet (sheet4): F2 357 -> 358, F4 76 -> 78
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 358
$ws4.Range("F4").Value = 78
